$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("17").Copy()
$ws.Rows("18").Insert()

Write-Host "Done insert"
